$wb = $excel.ActiveWorkbook

# --- ShopOrders sheet: add RevenueValue column ---
$ws1 = $wb.Worksheets.Item("ShopOrders")
$ws1.Range("O1").Value = "RevenueValue"
$ws1.Range("O1").Font.Bold = $true
$ws1.Range("O2").Value = 5
$ws1.Range("O3").Value = 2

# --- Selection / view state ---
$ws1.Activate()
$ws1.Application.ActiveWindow.ScrollColumn = 5
$ws1.Range("P10").Select()

$ws4 = $wb.Worksheets.Item("WorkCenterOpAllocations")
